$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (the source data keeps these as text).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.449.62"

$ws.Range("D3").Value = "1.824.51"

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "316.01"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").Value = "0.5169"
$ws.Range("E7").Value = "  +1.28%  "

$ws.Range("D8").Value = "0.3851"
$ws.Range("E8").Value = "  -1.54%  "

$ws.Range("D9").Value = "0.08301"
$ws.Range("E9").Value = "  +8.48%  "

$ws.Range("E10").Value = "  +1.25%  "

$ws.Range("D11").Value = "41.87"
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").Value = "6.375"
$ws.Range("E12").Value = "  +1.27%  "

$ws.Range("D13").Value = "21.09"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("D15").Value = "7.474"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").Value = "1.823.22"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").Value = "94.15"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").Value = "0.00001121"
$ws.Range("E18").Value = "  +3.49%  "

$ws.Range("D19").Value = "0.06634"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").Value = "6.052"
$ws.Range("E22").Value = "  -1.98%  "

$ws.Range("D23").Value = "28.483.49"

$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +2.67%  "

$ws.Range("D25").Value = "2.245"
$ws.Range("E25").Value = "  -0.62%  "

$ws.Range("D26").Value = "21.10"
$ws.Range("E26").Value = "  +2.24%  "

$ws.Range("D27").Value = "159.31"
$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("D28").Value = "2.035.02"
$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").Value = "2.403"
$ws.Range("E29").Value = "  +0.30%  "

$ws.Range("D30").Value = "125.94"
$ws.Range("E30").Value = "  +0.64%  "

$ws.Range("D31").Value = "0.1108"
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("D32").Value = "1.091"
$ws.Range("E32").Value = "  -3.10%  "

$ws.Range("D33").Value = "5.726"
$ws.Range("E33").Value = "  +0.76%  "

$ws.Range("D34").Value = "0.07520"
$ws.Range("E34").Value = "  +7.01%  "

$ws.Range("D35").Value = "3.685"
$ws.Range("E35").Value = "  +0.64%  "

$ws.Range("D36").Value = "0.2226"
$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("D37").Value = "0.02359"
$ws.Range("E37").Value = "  +1.44%  "

$ws.Range("D38").Value = "12.10"
$ws.Range("E38").Value = "  +7.66%  "

$ws.Range("D39").Value = "5.255"
$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("D40").Value = "8.760"
$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("D41").Value = "0.6390"
$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("D42").Value = "1.187"
$ws.Range("E42").Value = "  +0.31%  "

$ws.Range("D43").Value = "1.396"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "0.6193"
$ws.Range("E44").Value = "  +4.85%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.58"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("D46").Value = "3.796"
$ws.Range("E46").Value = "  +2.13%  "

$ws.Range("D47").Value = "127.35"
$ws.Range("E47").Value = "  +2.43%  "

$ws.Range("D48").Value = "2.003"
$ws.Range("E48").Value = "  +1.01%  "

$ws.Range("D49").Value = "1.204"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("D50").Value = "0.06961"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").Value = "1.083"
$ws.Range("E51").Value = "  +1.52%  "
